# Weekly update: insert a new weekly record for "Terminal Hortofrutícola
# Agro Chillán - Brócoli" at row 151, pushing the existing rows 151:231
# down to 152:232 (dimension grows from A1:R231 to A1:R232).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 151 (shifts 151:231 -> 152:232).
$ws.Rows(151).Insert()

# Populate the newly inserted row 151 with the new weekly observation.
$ws.Cells.Item(151, 1).Value = 7
$ws.Cells.Item(151, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(151, 3).Value = "Ñuble"
$ws.Cells.Item(151, 4).Value = 44609
$ws.Cells.Item(151, 5).Value = 16
$ws.Cells.Item(151, 6).Value = 100112023
$ws.Cells.Item(151, 7).Value = "Brócoli"
$ws.Cells.Item(151, 8).Value = "Sin especificar"
$ws.Cells.Item(151, 9).Value = "Primera"
$ws.Cells.Item(151, 10).Value = 200
$ws.Cells.Item(151, 11).Value = 700
$ws.Cells.Item(151, 12).Value = 750
$ws.Cells.Item(151, 13).Value = 725
$ws.Cells.Item(151, 14).Value = "$/unidad"
$ws.Cells.Item(151, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(151, 16).Value = 725
$ws.Cells.Item(151, 17).Value = 1
$ws.Cells.Item(151, 18).Value = "Hortaliza"
